# Daily refresh of the "剩余" (remaining) / "开始时间" (start date) tracker.
#
# Business logic (reverse-engineered from the data): each row tracks a
# countdown - D (总天, total days) minus elapsed days since F (开始时间,
# start date, stored as a yyyyMMdd integer) gives E (剩余, remaining
# days). This edit advances the "current day" by one (from 2025-10-09 to
# 2025-10-10): every row's elapsed-day count increases by 1, so E drops
# by 1 - unless that would bring the countdown to zero (or below), in
# which case the cycle restarts: F is reset to the new current day
# (20251010) and E is refilled back to D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = 20251010

function DateIntToOrdinal($d) {
    $s = [string]$d
    $y = [int]$s.Substring(0, 4)
    $m = [int]$s.Substring(4, 2)
    $dd = [int]$s.Substring(6, 2)
    $dt = Get-Date -Year $y -Month $m -Day $dd -Hour 0 -Minute 0 -Second 0
    return [math]::Floor($dt.ToOADate())
}

$todayOrdinal = DateIntToOrdinal $today

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $dCell = $ws.Cells.Item($row, 4)
    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $fStr = [string]([int]$fVal)

    # Row 36 carries a malformed start date (202510929) that can't be
    # parsed as yyyyMMdd - leave it untouched, same as the source diff.
    if ($fStr.Length -ne 8) {
        continue
    }

    $startOrdinal = DateIntToOrdinal ([int]$fVal)
    $elapsed = $todayOrdinal - $startOrdinal
    $newRemaining = [int]$dVal - $elapsed

    if ($newRemaining -le 0) {
        $eCell.Value2 = [int]$dVal
        $fCell.Value2 = $today
    } else {
        $eCell.Value2 = $newRemaining
    }
}
